$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 22 with change log entry
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat
$ws.Range("A22").Value = (Get-Date -Year 2022 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B22").Value = 2015
$ws.Range("C22").Value = "West Beach LOW"
$ws.Range("D22").Value = "Changed dates from 12 July 2014 to 14 July 2015"

$ws.Range("E22").Select() | Out-Null
